# FWMS - System UI changes.
# Applies the edits described by the commit diff against Open_Questions.xlsx.

$wb = $excel.ActiveWorkbook

$loginPage = $wb.Worksheets.Item("LoginPage")
$mapping   = $wb.Worksheets.Item("STOCK_SCREEN_MAPPING")

# --- Rename the second worksheet ---
$mapping.Name = "Sheet3"

# --- LoginPage sheet (sheet1.xml) ---
# D9 previously held "Correct as per Sreeni" - it's now a "Data Sent" status note.
$loginPage.Range("D9").Value = "Data Sent"
# F9 keeps its "Closed" value (unchanged content).
$loginPage.Range("F9").Value = "Closed"

# --- Sheet3 (formerly STOCK_SCREEN_MAPPING, sheet2.xml) ---
$sreeni = "SREENI - Ignore this field and no need to show this field"

$mapping.Range("E2").Value  = "I will send you code"
$mapping.Range("I2").Value  = $sreeni
$mapping.Range("I3").Value  = "this screen is master details - u should link to WMS_PRODUCT_ID field "

$mapping.Range("E9").Value  = "use for CUSIP"
$mapping.Range("E10").Value = $sreeni
$mapping.Range("E11").Value = $sreeni

$mapping.Range("E23").Value = $sreeni
$mapping.Range("E24").Value = "user for Divisor Days"

$mapping.Range("E29").Value = $sreeni
$mapping.Range("E32").Value = $sreeni
$mapping.Range("E34").Value = $sreeni
$mapping.Range("E35").Value = $sreeni
$mapping.Range("E36").Value = $sreeni
$mapping.Range("E37").Value = $sreeni
$mapping.Range("E41").Value = $sreeni
$mapping.Range("E44").Value = $sreeni
$mapping.Range("E47").Value = $sreeni
$mapping.Range("E49").Value = $sreeni
$mapping.Range("E50").Value = $sreeni

# New notes added under the "Bond Details" legend.
$mapping.Range("D55").Value = "use WMS_OTHER_CODE"
$mapping.Range("D56").Value = "use WMS_BOND_DIVISOR_DAYS_YEAR"

# Widen column E so the long "SREENI..." notes are fully visible (bestFit-style).
$mapping.Columns.Item(5).ColumnWidth = 51.28515625

# Move Sheet3's active selection from D62 to E3.
$mapping.Activate()
$mapping.Range("E3").Select()

# Move LoginPage's active selection from F10 to D10, and leave LoginPage as
# the selected/active tab (it was tabSelected="1" before the edit too).
$loginPage.Activate()
$loginPage.Range("D10").Select()
